# Update TPM-derived NATMI scores for the Fgf7-Fgfr3 LR-pair sheet.
# - Columns B/C (ligand/receptor symbol) are unchanged in content (Fgf7 / Fgfr3)
#   but a brand new cell type "Resolving-Mac" now exists in the cluster list, so
#   three additional rows (8-10) appear for the new sending cluster.
# - Existing numeric columns (I..T) get refreshed values reflecting the new TPM
#   input across the original 6 data rows (2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- refresh existing rows 2-7 (same cluster/gene labels, updated numbers) ----

$ws.Range("I2").Value = 0.9461861633339901
$ws.Range("J2").Value = 0.94618616333399
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.751166666666666
$ws.Range("N2").Value = 17.2535
$ws.Range("O2").Value = 0.7405222614421495
$ws.Range("P2").Value = 0.7405222614421495
$ws.Range("Q2").Value = 64.0720148151111
$ws.Range("R2").Value = 576.648133336
$ws.Range("S2").Value = 0.7006719174173573
$ws.Range("T2").Value = 0.7006719174173572

$ws.Range("I3").Value = 0.9461861633339901
$ws.Range("J3").Value = 0.94618616333399
$ws.Range("O3").Value = 0.07337387367415998
$ws.Range("P3").Value = 0.07337387367416
$ws.Range("S3").Value = 0.0694253440207063
$ws.Range("T3").Value = 0.0694253440207063

$ws.Range("I4").Value = 0.9461861633339901
$ws.Range("J4").Value = 0.94618616333399
$ws.Range("M4").Value = 1.445350666666667
$ws.Range("N4").Value = 4.336052
$ws.Range("O4").Value = 0.1861038648836906
$ws.Range("P4").Value = 0.1861038648836906
$ws.Range("Q4").Value = 16.10221624499911
$ws.Range("R4").Value = 144.919946204992
$ws.Range("S4").Value = 0.1760889018959265
$ws.Range("T4").Value = 0.1760889018959265

$ws.Range("G5").Value = 0.5338349999999999
$ws.Range("H5").Value = 1.601505
$ws.Range("I5").Value = 0.04533892403128163
$ws.Range("J5").Value = 0.04533892403128163
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.751166666666666
$ws.Range("N5").Value = 17.2535
$ws.Range("O5").Value = 0.7405222614421495
$ws.Range("P5").Value = 0.7405222614421495
$ws.Range("Q5").Value = 3.0701740575
$ws.Range("R5").Value = 27.6315665175
$ws.Range("S5").Value = 0.03357448255499849
$ws.Range("T5").Value = 0.03357448255499849

$ws.Range("G6").Value = 0.5338349999999999
$ws.Range("H6").Value = 1.601505
$ws.Range("I6").Value = 0.04533892403128163
$ws.Range("J6").Value = 0.04533892403128163
$ws.Range("O6").Value = 0.07337387367415998
$ws.Range("P6").Value = 0.07337387367416
$ws.Range("Q6").Value = 0.3042049850249999
$ws.Range("R6").Value = 2.737844865225
$ws.Range("S6").Value = 0.003326692484393594
$ws.Range("T6").Value = 0.003326692484393595

$ws.Range("G7").Value = 0.5338349999999999
$ws.Range("H7").Value = 1.601505
$ws.Range("I7").Value = 0.04533892403128163
$ws.Range("J7").Value = 0.04533892403128163
$ws.Range("M7").Value = 1.445350666666667
$ws.Range("N7").Value = 4.336052
$ws.Range("O7").Value = 0.1861038648836906
$ws.Range("P7").Value = 0.1861038648836906
$ws.Range("Q7").Value = 0.77157877314
$ws.Range("R7").Value = 6.944208958260001
$ws.Range("S7").Value = 0.008437748991889548
$ws.Range("T7").Value = 0.008437748991889548

# ---- new rows 8-10: sending cluster "Resolving-Mac" -> receptor cluster FAPs/MuSCs/ECs ----

$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("B8").Value = "Fgf7"
$ws.Range("C8").Value = "Fgfr3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.09978633333333332
$ws.Range("H8").Value = 0.299359
$ws.Range("I8").Value = 0.008474912634728231
$ws.Range("J8").Value = 0.008474912634728231
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.751166666666666
$ws.Range("N8").Value = 17.2535
$ws.Range("O8").Value = 0.7405222614421495
$ws.Range("P8").Value = 0.7405222614421495
$ws.Range("Q8").Value = 0.5738878340555554
$ws.Range("R8").Value = 5.1649905065
$ws.Range("S8").Value = 0.006275861469793595
$ws.Range("T8").Value = 0.006275861469793595

$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("B9").Value = "Fgf7"
$ws.Range("C9").Value = "Fgfr3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.09978633333333332
$ws.Range("H9").Value = 0.299359
$ws.Range("I9").Value = 0.008474912634728231
$ws.Range("J9").Value = 0.008474912634728231
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.5698483333333333
$ws.Range("N9").Value = 1.709545
$ws.Range("O9").Value = 0.07337387367415998
$ws.Range("P9").Value = 0.07337387367416
$ws.Range("Q9").Value = 0.05686307573944444
$ws.Range("R9").Value = 0.511767681655
$ws.Range("S9").Value = 0.0006218371690600916
$ws.Range("T9").Value = 0.0006218371690600917

$ws.Range("A10").Value = "Resolving-Mac"
$ws.Range("B10").Value = "Fgf7"
$ws.Range("C10").Value = "Fgfr3"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.09978633333333332
$ws.Range("H10").Value = 0.299359
$ws.Range("I10").Value = 0.008474912634728231
$ws.Range("J10").Value = 0.008474912634728231
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.445350666666667
$ws.Range("N10").Value = 4.336052
$ws.Range("O10").Value = 0.1861038648836906
$ws.Range("P10").Value = 0.1861038648836906
$ws.Range("Q10").Value = 0.1442262434075556
$ws.Range("R10").Value = 1.298036190668
$ws.Range("S10").Value = 0.001577213995874545
$ws.Range("T10").Value = 0.001577213995874545
